$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6653386454183267
$ws.Range("C2").Value = 0.6911764705882353
$ws.Range("D2").Value = 0.6573426573426573
$ws.Range("F2").Value = 0.7331300038124287

$ws.Range("B3").Value = 0.701195219123506
$ws.Range("C3").Value = 0.7126436781609196
$ws.Range("D3").Value = 0.7045454545454546
$ws.Range("E3").Value = 0.7209302325581395
$ws.Range("F3").Value = 0.7680772652179437
